$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 14
$ws.Range("B6").Value = 15
$ws.Range("B7").Value = 17
$ws.Range("B8").Value = 18
$ws.Range("B9").Value = 19
$ws.Range("B10").Value = 19
$ws.Range("B11").Value = 21
$ws.Range("B12").Value = 23
$ws.Range("B13").Value = 24
$ws.Range("B14").Value = 25
$ws.Range("B15").Value = 26
$ws.Range("B18").Value = 30
$ws.Range("B22").Value = 38
$ws.Range("B23").Value = 39
$ws.Range("B24").Value = 40
$ws.Range("B25").Value = 42
$ws.Range("B26").Value = 43
$ws.Range("B27").Value = 44
$ws.Range("B28").Value = 46
$ws.Range("B29").Value = 47
$ws.Range("B30").Value = 48
$ws.Range("B31").Value = 49
$ws.Range("B32").Value = 50
$ws.Range("B33").Value = 54
$ws.Range("B34").Value = 54
$ws.Range("B35").Value = 55
$ws.Range("B36").Value = 57
$ws.Range("B37").Value = 58
$ws.Range("B38").Value = 55
$ws.Range("B39").Value = 55
$ws.Range("B40").Value = 54
$ws.Range("B41").Value = 52
$ws.Range("B42").Value = 51
$ws.Range("B45").Value = 43
$ws.Range("B46").Value = 43
$ws.Range("B47").Value = 41
$ws.Range("B48").Value = 39
$ws.Range("B49").Value = 35
$ws.Range("B50").Value = 35
$ws.Range("B51").Value = 34
$ws.Range("B52").Value = 34
$ws.Range("B53").Value = 33
$ws.Range("B54").Value = 33
$ws.Range("B55").Value = 32
$ws.Range("B56").Value = 32
$ws.Range("B57").Value = 32
$ws.Range("B58").Value = 31
